# Update the Tuesday/Thursday schedule sheet with the new "Review Day" label
# and append the Lesson# -> Topic reference table beneath the calendar grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tuesday_Thursday_Schedule")

# Week 6 (row 7) already reads "Review Day" in column D; week 12 (row 12) gets
# the same label now that it's also being used as a review week.
$ws.Range("D12").Value = "Review Day"

# New lesson/topic lookup table appended below the schedule (rows 25-40).
# Column A (lesson numbers) was filled in for the first two rows before the
# column B topics were typed, matching the authoring order in the sheet.
$ws.Cells.Item(25, 1).Value = "Lesson 1"
$ws.Cells.Item(26, 1).Value = "Lesson 2"
$ws.Cells.Item(25, 2).Value = "Modeling Review"
$ws.Cells.Item(26, 2).Value = "Jupyter Review"

$lessonTopics = @(
    @("Lesson 3", "Network Models"),
    @("Lesson 4", "Shortest Path"),
    @("Lesson 5", "Modeling Functions in Jupyter"),
    @("Lesson 6", "Fixed Charge"),
    @("Lesson 7", "Set Covering"),
    @("Lesson 8", "Logical Constraints"),
    @("Lesson 9", "Python Review"),
    @("Lesson 10", "Min Spanning Tree"),
    @("Lesson 11", "TSP"),
    @("Lesson 12", "VRP"),
    @("Lesson 13", "Facility Location"),
    @("Lesson 14", "IP Formulations"),
    @("Lesson 15", "More IP Formulations"),
    @("Lesson 16", "Branch & Bound")
)

$startRow = 27
for ($i = 0; $i -lt $lessonTopics.Count; $i++) {
    $row = $startRow + $i
    $pair = $lessonTopics[$i]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}

# Reflect the author's final cursor position/selection.
$ws.Range("E11").Select()
